{"js": "// Apply the nanny-info text updates using the Word JavaScript API.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\n// Map of exact original paragraph text -> replacement text.\nconst replacements = new Map([\n  [\"This morning, Joey woke up at 12:00 AM.\", \"This morning, Joey woke up at 07:00 AM.\"],\n  [\"Her first nap should be at 12:00 AM.\", \"Her first nap should be at 09:00 AM.\"],\n  [\"For lunch today, we have A.\", \"For lunch today, we have Cherries.\"],\n  [\"For dinner today, we have B.\", \"For dinner today, we have Berries.\"],\n  [\"As a reminder C.\", \"As a reminder Don\\u2019t mess up.\"],\n]);\n\nfor (const para of paragraphs.items) {\n  const text = para.text;\n  if (replacements.has(text)) {\n    para.insertText(replacements.get(text), Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Apply the nanny-info text updates using the Word COM object model.\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute(\n        $findText,     # FindText\n        $false,        # MatchCase\n        $false,        # MatchWholeWord\n        $false,        # MatchWildcards\n        $false,        # MatchSoundsLike\n        $false,        # MatchAllWordForms\n        $true,         # Forward\n        1,             # Wrap = wdFindContinue\n        $false,        # Format\n        $replaceText,  # ReplaceWith\n        2              # Replace = wdReplaceAll\n    ) | Out-Null\n}\n\nReplace-Text \"This morning, Joey woke up at 12:00 AM.\" \"This morning, Joey woke up at 07:00 AM.\"\nReplace-Text \"Her first nap should be at 12:00 AM.\" \"Her first nap should be at 09:00 AM.\"\nReplace-Text \"For lunch today, we have A.\" \"For lunch today, we have Cherries.\"\nReplace-Text \"For dinner today, we have B.\" \"For dinner today, we have Berries.\"\nReplace-Text \"As a reminder C.\" \"As a reminder Don\u2019t mess up.\"\n"}
